# Add a "Continent" column to Sheet1 (inserted before Country), and tweak
# the selections on Sheet1 / Sheet2 to match the saved-file state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Map each country (Sheet1 column B after the insert) to its continent.
$continentByCountry = @{
    "United States"  = "North America"
    "Brazil"         = "South America"
    "India"          = "Asia"
    "Russia"         = "Transcontinental"
    "Peru"           = "South America"
    "Chile"          = "South America"
    "United Kingdom" = "Europe"
    "Mexico"         = "North America"
    "Spain"          = "Transcontinental"
    "Iran"           = "Asia"
    "Italy"          = "Europe"
    "Pakistan"       = "Asia"
    "Saudi Arabia"   = "Asia"
    "Turkey"         = "Transcontinental"
    "South Africa"   = "Africa"
    "Germany"        = "Europe"
    "France"         = "Europe"
    "Bangladesh"     = "Asia"
    "Colombia"       = "South America"
    "Canada"         = "North America"
}

# Insert a new column at A (shifts existing A:D -> B:E, widths travel with
# the cells automatically).
$ws1.Columns("A").Insert() | Out-Null

# Header for the new column, styled like the other header cells (copy the
# format from the adjacent header cell).
$ws1.Range("A1").Value = "Continent"
$ws1.Range("B1").Copy() | Out-Null
$ws1.Range("A1").PasteSpecial(-4122) | Out-Null

# Fill in the continent for each data row, looked up from the country that
# is now in column B.
for ($r = 2; $r -le 21; $r++) {
    $country = $ws1.Cells.Item($r, 2).Value2
    $ws1.Cells.Item($r, 1).Value = $continentByCountry[$country]
}

# Size the new column similarly to the data it holds (closest achievable
# match to the bestFit width Excel computed for "Transcontinental").
$ws1.Columns("A").ColumnWidth = 13.92

# Selections as left after the edit: Sheet2 had its whole column A selected,
# Sheet1 ends with B8 selected (and remains the active tab).
$ws2.Columns("A").Select() | Out-Null
$ws1.Range("B8").Select() | Out-Null
